# Updates cryptos list figures (price + 1h volume change) to the latest
# scrape, and fixes the Frax / EnergySwap row ordering (rows 47-48 swapped).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.310.10'
$ws.Range("E2").Value = '  +0.10%  '
$ws.Range("D3").Value = '1.690.78'
$ws.Range("E3").Value = '  +0.79%  '
$ws.Range("E4").Value = '  +0.23%  '
$ws.Range("D5").Value = "'" + '219.10'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.46%  '
$ws.Range("D6").Value = "'" + '0.5260'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.07%  '
$ws.Range("E7").Value = '  +0.27%  '
$ws.Range("D8").Value = "'" + '0.2705'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.41%  '
$ws.Range("D9").Value = "'" + '0.06444'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.18%  '
$ws.Range("D10").Value = "'" + '22.05'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.62%  '
$ws.Range("D11").Value = "'" + '0.07472'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.46%  '
$ws.Range("D12").Value = '1.706.19'
$ws.Range("E12").Value = '  +1.63%  '
$ws.Range("E13").Value = '  +0.13%  '
$ws.Range("D14").Value = "'" + '0.5857'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.52%  '
$ws.Range("D15").Value = "'" + '0.000008537'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.55%  '
$ws.Range("E16").Value = '  -1.03%  '
$ws.Range("D17").Value = '26.357.13'
$ws.Range("E17").Value = '  -0.04%  '
$ws.Range("D18").Value = "'" + '4.960'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.02%  '
$ws.Range("E19").Value = '  +0.16%  '
$ws.Range("D20").Value = "'" + '10.90'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.31%  '
$ws.Range("D21").Value = "'" + '189.79'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.35%  '
$ws.Range("D22").Value = "'" + '6.223'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.00%  '
$ws.Range("E23").Value = '  +0.23%  '
$ws.Range("D24").Value = "'" + '144.91'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.53%  '
$ws.Range("D25").Value = "'" + '7.666'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.39%  '
$ws.Range("D26").Value = "'" + '0.1235'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +5.01%  '
$ws.Range("E27").Value = '  +0.67%  '
$ws.Range("D28").Value = "'" + '0.06674'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +14.58%  '
$ws.Range("D29").Value = "'" + '1.354'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +5.70%  '
$ws.Range("E30").Value = '  +0.72%  '
$ws.Range("D31").Value = "'" + '3.589'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.87%  '
$ws.Range("D32").Value = "'" + '3.571'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.84%  '
$ws.Range("D33").Value = "'" + '1.672'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.11%  '
$ws.Range("D34").Value = "'" + '1.030'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.77%  '
$ws.Range("D35").Value = "'" + '0.6222'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +3.49%  '
$ws.Range("D36").Value = "'" + '2.393'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.52%  '
$ws.Range("E37").Value = '  +1.57%  '
$ws.Range("D38").Value = "'" + '6.356'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +5.21%  '
$ws.Range("D39").Value = "'" + '0.01626'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.44%  '
$ws.Range("D40").Value = '1.107.19'
$ws.Range("E40").Value = '  +2.33%  '
$ws.Range("D41").Value = "'" + '0.8841'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.69%  '
$ws.Range("E42").Value = '  +0.92%  '
$ws.Range("D43").Value = "'" + '100.94'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.98%  '
$ws.Range("D44").Value = '1.837.81'
$ws.Range("E44").Value = '  +0.58%  '
$ws.Range("D45").Value = "'" + '0.00000000116'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.87%  '
$ws.Range("D46").Value = "'" + '56.97'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.63%  '
$ws.Range("B47").Value = 'Frax'
$ws.Range("C47").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D47").Value = "'" + '1.012'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.54%  '
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").Value = "'" + '8.194'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.22%  '
$ws.Range("E49").Value = '  +1.48%  '
$ws.Range("E50").Value = '  +0.18%  '
$ws.Range("D51").Value = "'" + '6.061'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.89%  '
